$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Student ID column keeps text formatting (IDs are stored as text, not numbers)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A7").NumberFormat = "@"

# Row 2: Student ID 190874 -> 201252, Log Date 25/10/2025 -> 26/10/2025
$ws.Range("A2").Value = "201252"
$ws.Range("C2").Value = "26/10/2025"

# Row 3: Log Date 25/10/2025 -> 26/10/2025 (Student ID unchanged, 201253)
$ws.Range("C3").Value = "26/10/2025"

# Row 4: Student ID 201252 -> 200958, Log Date 25/10/2025 -> 26/10/2025
$ws.Range("A4").Value = "200958"
$ws.Range("C4").Value = "26/10/2025"

# Row 5: Student ID 201023 -> 201574, Log Date 25/10/2025 -> 26/10/2025
$ws.Range("A5").Value = "201574"
$ws.Range("C5").Value = "26/10/2025"

# Row 6: Student ID 201670 -> 201438, Log Date 25/10/2025 -> 26/10/2025
$ws.Range("A6").Value = "201438"
$ws.Range("C6").Value = "26/10/2025"

# Row 7: Student ID 190796 -> 211137, Log Date 25/10/2025 -> 26/10/2025
$ws.Range("A7").Value = "211137"
$ws.Range("C7").Value = "26/10/2025"

# Rows 8-13 no longer present in the updated data export; remove them entirely
$ws.Range("A8:F13").EntireRow.Delete()
